# Adding new test scripts for watch list
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- window / view bookkeeping (best effort) ---------------------------
$excel.ActiveWindow.WindowState = -4143
$excel.ActiveWindow.Width = 14310
$excel.ActiveWindow.Height = 4995

# Selection moves to D2:D22 and the frozen "topLeftCell" scroll position
# from the previous edit (C1) is cleared by simply re-selecting.
$ws.Range("D2:D22").Select()

# --- existing rows: Results column PASS -> SKIP (rows 2-19) ------------
$ws.Range("E2:E19").Value = "SKIP"

# --- three new test cases (rows 20-22) ----------------------------------
# Clone the formatting of the last existing data row (19) down into the
# three new rows first, so the new cells keep the same borders/wrap style
# without inventing any new cell-style entries.
$ws.Range("A19:E19").Copy()
$ws.Range("A20:E22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Rows.Item(20).RowHeight = 45
$ws.Rows.Item(21).RowHeight = 45
$ws.Rows.Item(22).RowHeight = 45

# TestCase_E19
$ws.Range("A20").Value = "TestCase_E19"
$ws.Range("B20").Value = "OPQA-288"
$ws.Range("C20").Value = "Verify that following fields are getting displayed for each article in the watchlist page:`na)Times cited`nb)Comments"
$ws.Range("D20").Value = "Y"
$ws.Range("E20").Value = "SKIP"

# TestCase_E20
$ws.Range("A21").Value = "TestCase_E20"
$ws.Range("B21").Value = "OPQA-290"
$ws.Range("C21").Value = "Verify that following fields are getting displayed for each article in the watchlist page:`na)Times cited`nb)Comments"
$ws.Range("D21").Value = "Y"
$ws.Range("E21").Value = "SKIP"

# TestCase_E21
$ws.Range("A22").Value = "TestCase_E21"
$ws.Range("B22").Value = "OPQA-291"
$ws.Range("C22").Value = "Verify that following fields are getting displayed for each post in the watchlist page:`na)Likes`nb)Comments"
$ws.Range("D22").Value = "Y"
$ws.Range("E22").Value = "PASS"
